# Manchester City fixtures workbook — update the attendance/number column (C)
# for all 20 fixture rows. The values in column C are plain text (shared
# strings) that happen to look like numbers, so we briefly mark the range as
# Text before assigning, then clear the number format back off again so the
# cells end up with no explicit style (matching how the sheet was originally
# authored) while keeping the values as text rather than numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    1  = "36"
    2  = "43"
    3  = "36"
    4  = "47"
    5  = "53"
    6  = "35"
    7  = "38"
    8  = "55"
    9  = "53"
    10 = "40"
    11 = "39"
    12 = "58"
    13 = "54"
    14 = "55"
    15 = "51"
    16 = "52"
    17 = "45"
    18 = "39"
    19 = "64"
    20 = "53"
}

$colRange = $ws.Range("C1:C20")
$colRange.NumberFormat = "@"

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 3).Value = $newValues[$row]
}

$colRange.ClearFormats()
